$d = $word.ActiveDocument

$replacements = @(
    @("89×87=7743", "33×55=1815"),
    @("20×86=1720", "33×92=3036"),
    @("69×94=6486", "75×60=4500"),
    @("52×31=1612", "43×33=1419"),
    @("62×56=3472", "28×91=2548"),
    @("60×76=4560", "85×19=1615"),
    @("26×53=1378", "65×81=5265"),
    @("11×98=1078", "16×74=1184"),
    @("45×80=3600", "69×36=2484"),
    @("96×26=2496", "21×81=1701"),
    @("44×84=3696", "37×80=2960"),
    @("70×56=3920", "98×23=2254"),
    @("76×77=5852", "52×29=1508"),
    @("16×45=720",  "31×23=713"),
    @("48×36=1728", "38×59=2242"),
    @("71×58=4118", "47×83=3901"),
    @("71×26=1846", "43×46=1978"),
    @("97×15=1455", "39×93=3627"),
    @("33×85=2805", "77×88=6776"),
    @("84×38=3192", "36×90=3240"),
    @("35×35=1225", "69×78=5382"),
    @("44×75=3300", "66×43=2838"),
    @("27×13=351",  "45×76=3420"),
    @("70×86=6020", "66×79=5214"),
    @("39×35=1365", "65×40=2600")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
